# Generate Report for Handback
# The localization-status report gets a freshly-generated handback entry for
# 9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3 on both the zh-cn and de-de sheets:
# a "Latest Target File" hyperlink, a "Latest Handback File" name, a
# "Latest Handback DateTime" stamp, and an "Error Detail" message explaining
# that the handed-back file is stale.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d126f93cb3233ed1c59a7c9bf41e229a2f66da1/e2e/9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c3464911cab506815d0d320d8687037da5fc3318/e2e/9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d126f93cb3233ed1c59a7c9bf41e229a2f66da1/e2e/9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.md."

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J8").Value = "9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.44e708b35d8f91b20c4bcf976d68672c0689e291.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-09-07 04:56:38"
$wsZh.Range("P8").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", "9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.md")

$wsZh.Columns.Item(16).ColumnWidth = 39.1666667

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J8").Value = "9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.44e708b35d8f91b20c4bcf976d68672c0689e291.de-de.xlf"
$wsDe.Range("K8").Value = "2016-09-07 04:56:47"
$wsDe.Range("P8").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", "9a0dc922-44b1-4e71-b24d-f3b1c7eeeba3.md")

$wsDe.Columns.Item(16).ColumnWidth = 39.1666667
